# edit.ps1 -- apply the KeyListenerDiagram slide-tree change described by the
# commit "Update documentations to reflect change".
#
# Net effect on ppt/slides/slide1.xml required by the diff:
#   * The shape "Rectangle 281" (id 282) is repositioned/resized:
#       off.x  4040813 -> 4114800   (off.y stays 2713934)
#       ext.cx  229325 -> 108049    (ext.cy stays 160062)
#     and it moves earlier in the z-order / shape tree: it used to sit right
#     after the 4th "Elbow Connector 63" (id 122), immediately before the
#     removed "Freeform 115" (id 136); it now sits right after the 2nd
#     "Elbow Connector 63" (id 85), immediately before the 3rd
#     "Elbow Connector 63" (id 97).
#   * The shape "Freeform 115" (id 136) -- the long connecting arrow that used
#     to live right after "Rectangle 281" at the end of the shape tree -- is
#     deleted outright.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

$EMU_PER_POINT = 12700

# --- Rectangle 281 (id 282): resize + move two slots back in the z-order ---
$rect = Get-ShapeById $s.Shapes 282

# New left / width (in points, since the Shape API works in points not EMU).
$rect.Left = 4114800 / $EMU_PER_POINT
$rect.Width = 108049 / $EMU_PER_POINT
# Top / Height are unchanged (2713934 / 160062 EMU respectively) so they are
# left untouched to avoid any needless floating point round-trip.

# Move the shape from after "Elbow Connector 63" (id 122) to right after
# "Elbow Connector 63" (id 85), i.e. two positions back, so the new order
# becomes: ... id=85, id=282(Rectangle 281), id=97, id=122 ...
$rect.ZOrder(3)  # msoSendBackward
$rect.ZOrder(3)  # msoSendBackward

# --- Freeform 115 (id 136): remove entirely ---
$freeform = Get-ShapeById $s.Shapes 136
$freeform.Delete()
